$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Holidays 2019")

$ws.Range("C2").Value = '30 42 груз сер'
$ws.Range("G2").Value = '30, 42, груз, сер'
$ws.Range("I2").Value = 45351
$ws.Range("G3").Value = '30, 42, груз, сер'
$ws.Range("I3").Value = 45351
$ws.Range("I4").Value = 45351
$ws.Range("I5").Value = 45351
$ws.Range("I6").Value = 45351
$ws.Range("I7").Value = 45351
$ws.Range("I8").Value = 45351
$ws.Range("I9").Value = 45351
$ws.Range("I10").Value = 45351
$ws.Range("C11").Value = 'Type груз сер LS-2'
$ws.Range("I11").Value = 45351
$ws.Range("I12").Value = 45351
$ws.Range("C13").Value = '202B C Type сер LS-2 H'
$ws.Range("G13").Value = 'Type, груз, сер, LS-2'
$ws.Range("I13").Value = 45351
$ws.Range("C14").Value = 'груз б/к сер'
$ws.Range("I14").Value = 45351
$ws.Range("G15").Value = '202B, C, Type, сер, LS-2, H'
$ws.Range("I15").Value = 45351
$ws.Range("G16").Value = '202B, C, Type, сер, LS-2, H'
$ws.Range("I16").Value = 45351
$ws.Range("C17").Value = '8 сх сер'
$ws.Range("G17").Value = '202B, C, Type, сер, LS-2, H'
$ws.Range("I17").Value = 45351
$ws.Range("G18").Value = 'груз, б/к, сер'
$ws.Range("I18").Value = 45351
$ws.Range("G19").Value = 'груз, б/к, сер'
$ws.Range("I19").Value = 45351
$ws.Range("G20").Value = 'груз, б/к, сер'
$ws.Range("I20").Value = 45351
$ws.Range("G21").Value = 'груз, б/к, сер'
$ws.Range("I21").Value = 45351
$ws.Range("I22").Value = 45351
$ws.Range("I23").Value = 45351
$ws.Range("I24").Value = 45351
$ws.Range("I25").Value = 45351
$ws.Range("I26").Value = 45351
$ws.Range("I27").Value = 45351
$ws.Range("I28").Value = 45351
$ws.Range("I29").Value = 45351
$ws.Range("I30").Value = 45351
